$d = $word.ActiveDocument

# --- 1. Strip the trailing "(đ)"/"(s)" grading-key markers from the eight
#        answer-option runs (Câu 1 a-d, Câu 2 a-d). Doing this via Find/Replace
#        on just the marker text (never retyping the quoted option text) avoids
#        Word's smart-quote autocorrect from turning the straight quotes that
#        are already in the document into curly ones. ---

$f1 = $d.Content
$f1.Find.Execute(' (đ)', $true, $false, $false, $false, $false, $true, 1, $false, '', 2)

$f2 = $d.Content
$f2.Find.Execute('(đ)', $true, $false, $false, $false, $false, $true, 1, $false, '', 2)

$f3 = $d.Content
$f3.Find.Execute('(s)', $true, $false, $false, $false, $false, $true, 1, $false, '', 2)

# --- 2. Add the "Đáp án: ĐĐĐS" summary paragraph right after item d) of Câu 1,
#        before the "Câu 2:" paragraph. ---

$r1 = $d.Content
$r1.Find.Execute('d) Lưu tệp sao lưu với tên _banthuam.sql và chọn "Import" để thực hiện việc sao lưu.', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
$p1 = $r1.Paragraphs(1)
$p1.Range.InsertParagraphAfter()

$r1b = $d.Content
$r1b.Find.Execute('d) Lưu tệp sao lưu với tên _banthuam.sql và chọn "Import" để thực hiện việc sao lưu.', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
$p1b = $r1b.Paragraphs(1)
$new1 = $p1b.Next()
$new1.Range.Text = "Đáp án: ĐĐĐS"

# --- 3. Add the "Đáp án: ĐSĐĐ" summary paragraph right after item d) of Câu 2,
#        at the very end of the document. ---

$r2 = $d.Content
$r2.Find.Execute('d) Sau khi truy vấn thành công, nhấn F5 để làm tươi lại danh sách CSDL.', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
$p2 = $r2.Paragraphs(1)
$p2.Range.InsertParagraphAfter()

$r2b = $d.Content
$r2b.Find.Execute('d) Sau khi truy vấn thành công, nhấn F5 để làm tươi lại danh sách CSDL.', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
$p2b = $r2b.Paragraphs(1)
$new2 = $p2b.Next()
$new2.Range.Text = "Đáp án: ĐSĐĐ"
